$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Hyperlinks.Delete()

$ws.Range("A1").Value = 8412989861
$ws.Range("B1").Value = "lokesh@84"

$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:lokesh@84", [Type]::Missing, [Type]::Missing, "lokesh@84")

$wb.Styles.Item("Hyperlink").Delete()

Write-Host "Underline before:" $ws.Range("B1").Font.Underline
$ws.Range("B1").Font.Underline = 0
Write-Host "Underline after:" $ws.Range("B1").Font.Underline
$ws.Range("B1").Font.Color = 16711680
